# Apply "contingencies with rene fine" edit:
#  - Extend header row 1 with two new columns P (14) and Q (15), using the
#    same style/format as the other header cells.
#  - For data rows 2-25:
#       column I: 1 -> 2
#       column K: 2 -> 1
#       column M: 1 -> 2
#       column O: 2 -> 1
#       add column P = 2
#       add column Q = 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Copy the formatting of the existing header cell O1 onto the two new
# header cells P1 and Q1, then set their values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2 through 25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P = 2 (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q = 2 (new)
}
